# Generate Report for Handback
# Adds a new handback entry (e9d423f6-7645-4077-a496-1c56187ed8bd) as row 4
# to the "Overview", "zh-cn" and "de-de" worksheets.

function Set-HyperlinkCell {
    param($ws, $cellRef, $url, $disp)
    # Hyperlinks.Add already applies Excel's built-in "Hyperlink" visual style
    # (underline + theme color) to the cell, matching the look of the other
    # hyperlink cells in the workbook (A2, A3, D2, D3, F2, F3, G2, G3, ...).
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $disp) | Out-Null
}

$wb = $excel.ActiveWorkbook

$uuid = "e9d423f6-7645-4077-a496-1c56187ed8bd"
$hash = "90e203c78433310c50aa295a45489f19911c1658"
$mdName = "e9d423f6-7645-4077-a496-1c56187ed8bd.md"
$inSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$urlMdMain = "https://github.com/OpenLocalizationTest/oltest/blob/90e203c78433310c50aa295a45489f19911c1658/e2e/e9d423f6-7645-4077-a496-1c56187ed8bd.md"
Set-HyperlinkCell $wsOverview "A4" $urlMdMain $mdName
$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnXlf = "e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.zh-cn.xlf"
$urlZhCnMd = "https://github.com/OpenLocalizationTest/oltest/blob/90e203c78433310c50aa295a45489f19911c1658/e2e/e9d423f6-7645-4077-a496-1c56187ed8bd.md"
$urlZhCnHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90e203c78433310c50aa295a45489f19911c1658/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.zh-cn.xlf"
$urlZhCnMd2 = "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/90e203c78433310c50aa295a45489f19911c1658/e2e/e9d423f6-7645-4077-a496-1c56187ed8bd.md"
$urlZhCnHandback = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90e203c78433310c50aa295a45489f19911c1658/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.zh-cn.xlf"

Set-HyperlinkCell $wsZhCn "A4" $urlZhCnMd $mdName
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $inSync
Set-HyperlinkCell $wsZhCn "D4" $urlZhCnHandoff $zhCnXlf
$wsZhCn.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E4").Value = "2016-03-23 06:45:14"
Set-HyperlinkCell $wsZhCn "F4" $urlZhCnMd2 $mdName
Set-HyperlinkCell $wsZhCn "G4" $urlZhCnHandback $zhCnXlf
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H4").Value = "2016-03-23 06:45:54"
$wsZhCn.Range("J4").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeXlf = "e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.de-de.xlf"
$urlDeDeMd = "https://github.com/OpenLocalizationTest/oltest/blob/90e203c78433310c50aa295a45489f19911c1658/e2e/e9d423f6-7645-4077-a496-1c56187ed8bd.md"
$urlDeDeHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90e203c78433310c50aa295a45489f19911c1658/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.de-de.xlf"
$urlDeDeMd2 = "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/90e203c78433310c50aa295a45489f19911c1658/e2e/e9d423f6-7645-4077-a496-1c56187ed8bd.md"
$urlDeDeHandback = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90e203c78433310c50aa295a45489f19911c1658/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e9d423f6-7645-4077-a496-1c56187ed8bd.90e203c78433310c50aa295a45489f19911c1658.de-de.xlf"

Set-HyperlinkCell $wsDeDe "A4" $urlDeDeMd $mdName
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $inSync
Set-HyperlinkCell $wsDeDe "D4" $urlDeDeHandoff $deDeXlf
$wsDeDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("E4").Value = "2016-03-23 06:45:23"
Set-HyperlinkCell $wsDeDe "F4" $urlDeDeMd2 $mdName
Set-HyperlinkCell $wsDeDe "G4" $urlDeDeHandback $deDeXlf
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H4").Value = "2016-03-23 06:46:11"
$wsDeDe.Range("J4").Value = "Include"
